$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Ministero dell'Istruzione, dell'Universita' e della Ricerca -> Ministero dell'Istruzione, dell'Universita' e della Ricerca
$ws.Range("C2").Value = 8552

# Row 3: Unione Italiana delle Camere di Commercio Industria, Artigianato e Agricoltura -> Unione Italiana delle Camere di Commercio Industria, Artigianato e Agricoltura
$ws.Range("C3").Value = 1069

# Row 4: Poste Italiane S.P.A. -> Poste Italiane S.P.A.
$ws.Range("C4").Value = 779

# Row 5: E-Fil S.r.l. -> E-Fil S.r.l.
$ws.Range("C5").Value = 543

# Row 6: Regione Lombardia -> Regione Lombardia
$ws.Range("C6").Value = 471

# Row 7: Regione del Veneto -> Regione del Veneto
$ws.Range("C7").Value = 418

# Row 8: CREDEMTEL SpA -> Maggioli SPA
$ws.Range("A8").Value = "Maggioli SPA"
$ws.Range("B8").NumberFormat = "@"
$ws.Range("B8").Value = "06188330150"
$ws.Range("C8").Value = 405

# Row 9: Maggioli SPA -> CREDEMTEL SpA
$ws.Range("A9").Value = "CREDEMTEL SpA"
$ws.Range("B9").NumberFormat = "@"
$ws.Range("B9").Value = "01378570350"
$ws.Range("C9").Value = 405

# Row 11: Provincia Autonoma di Trento -> Banca Popolare di Sondrio, Società Cooperativa per Azioni
$ws.Range("A11").Value = "Banca Popolare di Sondrio, Società Cooperativa per Azioni"
$ws.Range("B11").NumberFormat = "@"
$ws.Range("B11").Value = "00053810149"
$ws.Range("C11").Value = 339

# Row 12: Banca Popolare di Sondrio, Società Cooperativa per Azioni -> Provincia Autonoma di Trento
$ws.Range("A12").Value = "Provincia Autonoma di Trento"
$ws.Range("B12").NumberFormat = "@"
$ws.Range("B12").Value = "00337460224"
$ws.Range("C12").Value = 334

# Row 13: Dedagroup Public Services S.R.L. -> Dedagroup Public Services S.R.L.
$ws.Range("C13").Value = 242

# Row 14: Progetti e Soluzioni SPA -> Progetti e Soluzioni SPA
$ws.Range("C14").Value = 188

# Row 16: P.A. Digitale spa -> APKAPPA S.R.L.
$ws.Range("A16").Value = "APKAPPA S.R.L."
$ws.Range("B16").NumberFormat = "@"
$ws.Range("B16").Value = "08543640158"
$ws.Range("C16").Value = 175

# Row 17: APKAPPA S.R.L. -> P.A. Digitale spa
$ws.Range("A17").Value = "P.A. Digitale spa"
$ws.Range("B17").NumberFormat = "@"
$ws.Range("B17").Value = "06628860964"
$ws.Range("C17").Value = 161

# Row 18: Regione Marche -> Regione Piemonte
$ws.Range("A18").Value = "Regione Piemonte"
$ws.Range("B18").NumberFormat = "@"
$ws.Range("B18").Value = "80087670016"
$ws.Range("C18").Value = 160

# Row 19: ADVANCED SYSTEMS srl -> Regione Marche
$ws.Range("A19").Value = "Regione Marche"
$ws.Range("B19").NumberFormat = "@"
$ws.Range("B19").Value = "80008630420"
$ws.Range("C19").Value = 157

# Row 20: Regione Piemonte -> Advanced Systems S.p.A.
$ws.Range("A20").Value = "Advanced Systems S.p.A."
$ws.Range("B20").NumberFormat = "@"
$ws.Range("B20").Value = "03383350638"
$ws.Range("C20").Value = 149

# Row 21: Regione Puglia -> Regione Autonoma della Sardegna
$ws.Range("A21").Value = "Regione Autonoma della Sardegna"
$ws.Range("B21").NumberFormat = "@"
$ws.Range("B21").Value = "80002870923"
$ws.Range("C21").Value = 129

# Row 22: Regione Basilicata -> Regione Puglia
$ws.Range("A22").Value = "Regione Puglia"
$ws.Range("B22").NumberFormat = "@"
$ws.Range("B22").Value = "80017210727"
$ws.Range("C22").Value = 122

# Row 23: Regione Toscana -> Regione Basilicata
$ws.Range("A23").Value = "Regione Basilicata"
$ws.Range("B23").NumberFormat = "@"
$ws.Range("B23").Value = "80002950766"
$ws.Range("C23").Value = 106

# Row 24: ANCITEL -> Regione Toscana
$ws.Range("A24").Value = "Regione Toscana"
$ws.Range("B24").NumberFormat = "@"
$ws.Range("B24").Value = "01386030488"
$ws.Range("C24").Value = 104

# Row 25: Next Step Solution -> ANCITEL
$ws.Range("A25").Value = "ANCITEL"
$ws.Range("B25").NumberFormat = "@"
$ws.Range("B25").Value = "07196850585"
$ws.Range("C25").Value = 90

# Row 26: Regione Autonoma Friuli-Venezia Giulia -> Next Step Solution
$ws.Range("A26").Value = "Next Step Solution"
$ws.Range("B26").NumberFormat = "@"
$ws.Range("B26").Value = "02554480349"
$ws.Range("C26").Value = 89

# Row 27: Regione Autonoma della Sardegna -> Regione Autonoma Friuli-Venezia Giulia
$ws.Range("A27").Value = "Regione Autonoma Friuli-Venezia Giulia"
$ws.Range("B27").NumberFormat = "@"
$ws.Range("B27").Value = "80014930327"
$ws.Range("C27").Value = 85

# Row 28: Intesa Sanpaolo SPA -> PMPay s.r.l.
$ws.Range("A28").Value = "PMPay s.r.l."
$ws.Range("B28").NumberFormat = "@"
$ws.Range("B28").Value = "08747230962"
$ws.Range("C28").Value = 67

# Row 29: PMPay s.r.l. -> Intesa Sanpaolo SPA
$ws.Range("A29").Value = "Intesa Sanpaolo SPA"
$ws.Range("B29").NumberFormat = "@"
$ws.Range("B29").Value = "00799960158"
$ws.Range("C29").Value = 64

# Row 30: NORDCOM -> Siscom SPA
$ws.Range("A30").Value = "Siscom SPA"
$ws.Range("B30").NumberFormat = "@"
$ws.Range("B30").Value = "01778000040"
$ws.Range("C30").Value = 62

# Row 31: Regione Umbria -> ROMA CAPITALE
$ws.Range("A31").Value = "ROMA CAPITALE"
$ws.Range("B31").NumberFormat = "@"
$ws.Range("B31").Value = "02438750586"
$ws.Range("C31").Value = 58

# Row 32: Bluenext S.r.l. -> NORDCOM
$ws.Range("A32").Value = "NORDCOM"
$ws.Range("B32").NumberFormat = "@"
$ws.Range("B32").Value = "13384100155"
$ws.Range("C32").Value = 56

# Row 33: ROMA CAPITALE -> Regione Umbria
$ws.Range("A33").Value = "Regione Umbria"
$ws.Range("B33").NumberFormat = "@"
$ws.Range("B33").Value = "80000130544"
$ws.Range("C33").Value = 52

# Row 34: Italriscossioni Società Italiana di Fiscalità Locale S.r.l. -> Bluenext S.r.l.
$ws.Range("A34").Value = "Bluenext S.r.l."
$ws.Range("B34").NumberFormat = "@"
$ws.Range("B34").Value = "04228480408"
$ws.Range("C34").Value = 50

# Row 35: CINECA consorzio universitario -> Italriscossioni Società Italiana di Fiscalità Locale S.r.l.
$ws.Range("A35").Value = "Italriscossioni Società Italiana di Fiscalità Locale S.r.l."
$ws.Range("B35").NumberFormat = "@"
$ws.Range("B35").Value = "06092371001"
$ws.Range("C35").Value = 48

# Row 36: Regione Autonoma Valle D'Aosta -> Regione Autonoma Valle D'Aosta
$ws.Range("C36").Value = 48

# Row 37: Consorzio I.T. Srl -> CINECA consorzio universitario
$ws.Range("A37").Value = "CINECA consorzio universitario"
$ws.Range("B37").NumberFormat = "@"
$ws.Range("B37").Value = "00317740371"
$ws.Range("C37").Value = 43

# Row 38: UNIMATICA S.P.A -> Consorzio I.T. Srl
$ws.Range("A38").Value = "Consorzio I.T. Srl"
$ws.Range("B38").NumberFormat = "@"
$ws.Range("B38").Value = "01321400192"
$ws.Range("C38").Value = 40

# Row 39: Siscom SPA -> UNIMATICA S.P.A
$ws.Range("A39").Value = "UNIMATICA S.P.A"
$ws.Range("B39").NumberFormat = "@"
$ws.Range("B39").Value = "02098391200"
$ws.Range("C39").Value = 39

# Row 43: Regione Liguria -> Regione Liguria
$ws.Range("C43").Value = 26

# Row 44: Novares Spa -> Novares Spa
$ws.Range("C44").Value = 23

# Row 45: ANDREANI TRIBUTI srl -> Numera Sistemi e Informatica SpA
$ws.Range("A45").Value = "Numera Sistemi e Informatica SpA"
$ws.Range("B45").NumberFormat = "@"
$ws.Range("B45").Value = "01265230902"
$ws.Range("C45").Value = 21

# Row 48: Numera Sistemi e Informatica SpA -> Regione Lazio
$ws.Range("A48").Value = "Regione Lazio"
$ws.Range("B48").NumberFormat = "@"
$ws.Range("B48").Value = "80143490581"

# Row 49: Regione Lazio -> ANDREANI TRIBUTI srl
$ws.Range("A49").Value = "ANDREANI TRIBUTI srl"
$ws.Range("B49").NumberFormat = "@"
$ws.Range("B49").Value = "01412920439"

# Row 51: Servizi Locali SpA -> Servizi Locali SpA
$ws.Range("C51").Value = 15

# Row 52: Crédit Agricole Group Solutions Società Consortile per azioni -> Si.Form Consulting srl
$ws.Range("A52").Value = "Si.Form Consulting srl"
$ws.Range("B52").NumberFormat = "@"
$ws.Range("B52").Value = "03943960827"
$ws.Range("C52").Value = 13

# Row 53: Si.Form Consulting srl -> Crédit Agricole Group Solutions Società Consortile per azioni
$ws.Range("A53").Value = "Crédit Agricole Group Solutions Società Consortile per azioni"
$ws.Range("B53").NumberFormat = "@"
$ws.Range("B53").Value = "02771790348"
$ws.Range("C53").Value = 12

# Row 58: Comune di Catania -> Comune di Catania
$ws.Range("C58").Value = 7

# Row 60: ISWEB S.p.A. -> Phoenix IT Solutions S.r.L
$ws.Range("A60").Value = "Phoenix IT Solutions S.r.L"
$ws.Range("B60").NumberFormat = "@"
$ws.Range("B60").Value = "07623321218"

# Row 62: Phoenix IT Solutions S.r.L -> CityPoste Payment Digital S.r.l.
$ws.Range("A62").Value = "CityPoste Payment Digital S.r.l."
$ws.Range("B62").NumberFormat = "@"
$ws.Range("B62").Value = "02003750672"
$ws.Range("C62").Value = 3

# Row 63: ICCREA Banca SpA -> Argentea S.r.l.
$ws.Range("A63").Value = "Argentea S.r.l."
$ws.Range("B63").NumberFormat = "@"
$ws.Range("B63").Value = "02260390220"
$ws.Range("C63").Value = 3

# Row 64: CityPoste Payment Digital S.r.l. -> ISWEB S.p.A.
$ws.Range("A64").Value = "ISWEB S.p.A."
$ws.Range("B64").NumberFormat = "@"
$ws.Range("B64").Value = "01722270665"
$ws.Range("C64").Value = 3

# Row 65: Agenzia Italiana del Farmaco - AIFA -> ICCREA Banca SpA
$ws.Range("A65").Value = "ICCREA Banca SpA"
$ws.Range("B65").NumberFormat = "@"
$ws.Range("B65").Value = "04774801007"
$ws.Range("C65").Value = 2

# Row 66: I.C.A. - Imposte Comunali Affini – s.r.l. -> Ministero dello Sviluppo Economico
$ws.Range("A66").Value = "Ministero dello Sviluppo Economico"
$ws.Range("B66").NumberFormat = "@"
$ws.Range("B66").Value = "80230390587"

# Row 67: MegASP S.r.l. -> BANCA MONTE DEI PASCHI DI SIENA
$ws.Range("A67").Value = "BANCA MONTE DEI PASCHI DI SIENA"
$ws.Range("B67").NumberFormat = "@"
$ws.Range("B67").Value = "00884060526"

# Row 68: Argentea S.r.l. -> Softline srl
$ws.Range("A68").Value = "Softline srl"
$ws.Range("B68").NumberFormat = "@"
$ws.Range("B68").Value = "12299030150"

# Row 69: Ministero dello Sviluppo Economico -> Banco BPM Società per Azioni
$ws.Range("A69").Value = "Banco BPM Società per Azioni"
$ws.Range("B69").NumberFormat = "@"
$ws.Range("B69").Value = "09722490969"

# Row 70: BANCA MONTE DEI PASCHI DI SIENA -> I.C.A. - Imposte Comunali Affini – s.r.l.
$ws.Range("A70").Value = "I.C.A. - Imposte Comunali Affini – s.r.l."
$ws.Range("B70").NumberFormat = "@"
$ws.Range("B70").Value = "02478610583"

# Row 71: Softline srl -> ARGO SOFTWARE SRL
$ws.Range("A71").Value = "ARGO SOFTWARE SRL"
$ws.Range("B71").NumberFormat = "@"
$ws.Range("B71").Value = "00838520880"

# Row 72: San Marco SPA -> Engineering Ingegneria Informatica SpA
$ws.Range("A72").Value = "Engineering Ingegneria Informatica SpA"
$ws.Range("B72").NumberFormat = "@"
$ws.Range("B72").Value = "00967720285"

# Row 73: Banco BPM Società per Azioni -> San Marco SPA
$ws.Range("A73").Value = "San Marco SPA"
$ws.Range("B73").NumberFormat = "@"
$ws.Range("B73").Value = "04142440728"

# Row 74: ARGO SOFTWARE SRL -> Noviservice srl
$ws.Range("A74").Value = "Noviservice srl"
$ws.Range("B74").NumberFormat = "@"
$ws.Range("B74").Value = "02789990922"

# Row 75: Engineering Ingegneria Informatica SpA -> Agenzia Italiana del Farmaco - AIFA
$ws.Range("A75").Value = "Agenzia Italiana del Farmaco - AIFA"
$ws.Range("B75").NumberFormat = "@"
$ws.Range("B75").Value = "97345810580"

# Row 76: Noviservice srl -> MegASP S.r.l.
$ws.Range("A76").Value = "MegASP S.r.l."
$ws.Range("B76").NumberFormat = "@"
$ws.Range("B76").Value = "09898030151"

